$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / safe-string cell updates (coin names, links, and price strings
# that Excel's numeric parser won't mistake for a number, e.g. multi-dot values)
$textUpdates = @(
    @{ Ref = 'D2'; Value = '21.181.03' }
    @{ Ref = 'E2'; Value = '  -3.98%  ' }
    @{ Ref = 'D3'; Value = '1.516.85' }
    @{ Ref = 'E3'; Value = '  -2.39%  ' }
    @{ Ref = 'E4'; Value = '  +0.70%  ' }
    @{ Ref = 'E5'; Value = '  +0.59%  ' }
    @{ Ref = 'E6'; Value = '  -1.20%  ' }
    @{ Ref = 'E7'; Value = '  -1.17%  ' }
    @{ Ref = 'E8'; Value = '  -2.15%  ' }
    @{ Ref = 'E9'; Value = '  -2.93%  ' }
    @{ Ref = 'E10'; Value = '  -2.81%  ' }
    @{ Ref = 'E11'; Value = '  -1.91%  ' }
    @{ Ref = 'E12'; Value = '  +0.83%  ' }
    @{ Ref = 'E13'; Value = '  -0.05%  ' }
    @{ Ref = 'E14'; Value = '  -3.76%  ' }
    @{ Ref = 'D15'; Value = '1.523.01' }
    @{ Ref = 'E15'; Value = '  -1.43%  ' }
    @{ Ref = 'E16'; Value = '  -3.56%  ' }
    @{ Ref = 'E17'; Value = '  -4.12%  ' }
    @{ Ref = 'E18'; Value = '  -0.20%  ' }
    @{ Ref = 'E19'; Value = '  -1.23%  ' }
    @{ Ref = 'E20'; Value = '  +0.63%  ' }
    @{ Ref = 'E21'; Value = '  -4.06%  ' }
    @{ Ref = 'E22'; Value = '  -1.47%  ' }
    @{ Ref = 'E23'; Value = '  -4.78%  ' }
    @{ Ref = 'E24'; Value = '  -0.88%  ' }
    @{ Ref = 'D25'; Value = '21.190.11' }
    @{ Ref = 'E25'; Value = '  -3.98%  ' }
    @{ Ref = 'E26'; Value = '  -1.97%  ' }
    @{ Ref = 'E27'; Value = '  -0.69%  ' }
    @{ Ref = 'E28'; Value = '  -1.99%  ' }
    @{ Ref = 'E29'; Value = '  -1.30%  ' }
    @{ Ref = 'D30'; Value = '1.693.20' }
    @{ Ref = 'E30'; Value = '  -1.77%  ' }
    @{ Ref = 'E31'; Value = '  -2.48%  ' }
    @{ Ref = 'E32'; Value = '  +3.40%  ' }
    @{ Ref = 'E33'; Value = '  -4.41%  ' }
    @{ Ref = 'E34'; Value = '  -3.55%  ' }
    @{ Ref = 'E35'; Value = '  -6.26%  ' }
    @{ Ref = 'E36'; Value = '  -0.23%  ' }
    @{ Ref = 'E37'; Value = '  -7.53%  ' }
    @{ Ref = 'E38'; Value = '  +6.53%  ' }
    @{ Ref = 'E39'; Value = '  -3.55%  ' }
    @{ Ref = 'E40'; Value = '  -4.88%  ' }
    @{ Ref = 'B41'; Value = 'Frax' }
    @{ Ref = 'C41'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax' }
    @{ Ref = 'E41'; Value = '  +0.54%  ' }
    @{ Ref = 'B42'; Value = 'Algorand' }
    @{ Ref = 'C42'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Ref = 'E42'; Value = '  -2.03%  ' }
    @{ Ref = 'E43'; Value = '  -3.80%  ' }
    @{ Ref = 'E44'; Value = '  -2.25%  ' }
    @{ Ref = 'E45'; Value = '  +0.58%  ' }
    @{ Ref = 'E46'; Value = '  -1.25%  ' }
    @{ Ref = 'E47'; Value = '  -1.64%  ' }
    @{ Ref = 'B48'; Value = 'EOS' }
    @{ Ref = 'C48'; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos' }
    @{ Ref = 'E48'; Value = '  +1.15%  ' }
    @{ Ref = 'B49'; Value = 'NEARProtocol' }
    @{ Ref = 'C49'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Ref = 'E49'; Value = '  -1.83%  ' }
    @{ Ref = 'E50'; Value = '  -2.14%  ' }
    @{ Ref = 'E51'; Value = '  -3.50%  ' }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Ref).Value = $u.Value
}

# Numeric-looking price strings: Excel's COM layer auto-converts plain numeric
# text to a real number (and silently normalises e.g. trailing zeros) when you
# assign via .Value. Force text storage by switching the cell to a text number
# format just for the assignment, then clear formatting again so the cell keeps
# its original (unstyled) appearance, matching the source workbook.
$numericTextUpdates = @(
    @{ Ref = 'D4'; Value = '1.007' }
    @{ Ref = 'D5'; Value = '1.007' }
    @{ Ref = 'D6'; Value = '286.73' }
    @{ Ref = 'D7'; Value = '0.3888' }
    @{ Ref = 'D8'; Value = '0.3144' }
    @{ Ref = 'D9'; Value = '42.37' }
    @{ Ref = 'D10'; Value = '0.07035' }
    @{ Ref = 'D11'; Value = '1.052' }
    @{ Ref = 'D13'; Value = '5.650' }
    @{ Ref = 'D16'; Value = '6.376' }
    @{ Ref = 'D18'; Value = '0.06582' }
    @{ Ref = 'D19'; Value = '82.26' }
    @{ Ref = 'D20'; Value = '1.007' }
    @{ Ref = 'D21'; Value = '6.023' }
    @{ Ref = 'D22'; Value = '15.21' }
    @{ Ref = 'D23'; Value = '10.71' }
    @{ Ref = 'D24'; Value = '2.354' }
    @{ Ref = 'D26'; Value = '2.375' }
    @{ Ref = 'D27'; Value = '147.51' }
    @{ Ref = 'D28'; Value = '18.13' }
    @{ Ref = 'D29'; Value = '4.816' }
    @{ Ref = 'D31'; Value = '115.43' }
    @{ Ref = 'D32'; Value = '5.983' }
    @{ Ref = 'D33'; Value = '0.9500' }
    @{ Ref = 'D34'; Value = '0.07990' }
    @{ Ref = 'D35'; Value = '8.441' }
    @{ Ref = 'D36'; Value = '5.088' }
    @{ Ref = 'D37'; Value = '1.485' }
    @{ Ref = 'D38'; Value = '11.34' }
    @{ Ref = 'D39'; Value = '0.05837' }
    @{ Ref = 'D40'; Value = '0.02147' }
    @{ Ref = 'D41'; Value = '1.006' }
    @{ Ref = 'D42'; Value = '0.1993' }
    @{ Ref = 'D43'; Value = '1.165' }
    @{ Ref = 'D44'; Value = '0.5666' }
    @{ Ref = 'D45'; Value = '12.95' }
    @{ Ref = 'D46'; Value = '3.697' }
    @{ Ref = 'D47'; Value = '0.5480' }
    @{ Ref = 'D48'; Value = '1.142' }
    @{ Ref = 'D49'; Value = '1.856' }
    @{ Ref = 'D50'; Value = '115.02' }
    @{ Ref = 'D51'; Value = '0.06577' }
)

foreach ($u in $numericTextUpdates) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

